$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 1915
$ws.Range("H3").Value = 8348
$ws.Range("I3").Value = 2028
$ws.Range("G4").Value = 1426
$ws.Range("I4").Value = 508
$ws.Range("I5").Value = 176
$ws.Range("I6").Value = 2452
$ws.Range("G7").Value = 24645
$ws.Range("H7").Value = 25970
$ws.Range("I7").Value = 7079

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I2").Value = 18
$ws.Range("I7").Value = 72

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("I3").Value = 7
$ws.Range("I7").Value = 23

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 73
$ws.Range("I3").Value = 69
$ws.Range("I7").Value = 231

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I6").Value = 36
$ws.Range("I7").Value = 127

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 95
$ws.Range("I6").Value = 92
$ws.Range("I7").Value = 270

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I6").Value = 25
$ws.Range("I7").Value = 69

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I6").Value = 61
$ws.Range("I7").Value = 165

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("H8").Value = 1698
$ws.Range("I8").Value = 441
$ws.Range("I11").Value = 123
$ws.Range("I18").Value = 58
$ws.Range("I19").Value = 207
$ws.Range("G20").Value = 620
$ws.Range("I20").Value = 190
$ws.Range("I22").Value = 20
$ws.Range("I27").Value = 66
$ws.Range("I29").Value = 461
$ws.Range("I30").Value = 23
$ws.Range("I31").Value = 69
$ws.Range("I33").Value = 335
$ws.Range("I34").Value = 30
$ws.Range("I37").Value = 231
$ws.Range("I40").Value = 11
$ws.Range("I42").Value = 238
$ws.Range("I44").Value = 55
$ws.Range("I52").Value = 146
$ws.Range("I54").Value = 156
$ws.Range("G63").Value = 189
$ws.Range("I63").Value = 33
$ws.Range("I65").Value = 165
$ws.Range("I67").Value = 270
$ws.Range("I74").Value = 17
$ws.Range("I85").Value = 335
$ws.Range("I86").Value = 44
$ws.Range("I88").Value = 55
$ws.Range("I89").Value = 72
$ws.Range("I92").Value = 21
$ws.Range("I93").Value = 35
$ws.Range("I94").Value = 57
$ws.Range("I99").Value = 127
$ws.Range("G101").Value = 24645
$ws.Range("H101").Value = 25970
$ws.Range("I101").Value = 7079

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I2").Value = 83
$ws.Range("I3").Value = 116
$ws.Range("I7").Value = 335

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I3").Value = 33
$ws.Range("I7").Value = 156

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I3").Value = 150
$ws.Range("I6").Value = 145
$ws.Range("I7").Value = 461

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I2").Value = 86
$ws.Range("I3").Value = 51
$ws.Range("I7").Value = 207

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("I3").Value = 11
$ws.Range("I7").Value = 55

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I6").Value = 96
$ws.Range("I7").Value = 335

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I3").Value = 83
$ws.Range("I6").Value = 63
$ws.Range("I7").Value = 238

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I3").Value = 58
$ws.Range("G4").Value = 26
$ws.Range("I6").Value = 70
$ws.Range("G7").Value = 620
$ws.Range("I7").Value = 190

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I2").Value = 13
$ws.Range("I7").Value = 58

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("I3").Value = 11
$ws.Range("I7").Value = 35

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I6").Value = 28
$ws.Range("I7").Value = 146

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("I6").Value = 10
$ws.Range("I7").Value = 30

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I3").Value = 7
$ws.Range("I6").Value = 33
$ws.Range("I7").Value = 57

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I5").Value = 8
$ws.Range("I7").Value = 123

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("I2").Value = 7
$ws.Range("I7").Value = 21

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I6").Value = 21
$ws.Range("I7").Value = 55

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 136
$ws.Range("H3").Value = 587
$ws.Range("I3").Value = 119
$ws.Range("I6").Value = 147
$ws.Range("H7").Value = 1698
$ws.Range("I7").Value = 441

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I6").Value = 32
$ws.Range("I7").Value = 66

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I4").Value = 30
$ws.Range("I7").Value = 44

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("I6").Value = 10
$ws.Range("I7").Value = 20

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("I3").Value = 4
$ws.Range("I7").Value = 11

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("I6").Value = 12
$ws.Range("I7").Value = 17
